# Update "想去人数" (number of people interested) figures that changed
# between the two scrapes, on both the "展览" (Exhibition) sheet and the
# "全部类型" (All types) sheet, which both list the same events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> column F, rows keyed by row number
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 258
$wsExhibition.Range("F6").Value  = 10090
$wsExhibition.Range("F8").Value  = 913
$wsExhibition.Range("F9").Value  = 1255
$wsExhibition.Range("F10").Value = 6068
$wsExhibition.Range("F12").Value = 413
$wsExhibition.Range("F13").Value = 186
$wsExhibition.Range("F15").Value = 3100
$wsExhibition.Range("F18").Value = 597
$wsExhibition.Range("F23").Value = 1540

# Sheet "全部类型" -> column F, rows keyed by row number (shifted by +1
# relative to "展览" because it includes an extra performance event)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 258
$wsAll.Range("F7").Value  = 10090
$wsAll.Range("F9").Value  = 913
$wsAll.Range("F10").Value = 1255
$wsAll.Range("F11").Value = 6068
$wsAll.Range("F13").Value = 413
$wsAll.Range("F14").Value = 186
$wsAll.Range("F16").Value = 3100
$wsAll.Range("F19").Value = 597
$wsAll.Range("F24").Value = 1540
